# Update countries & provincias Spain
# Refreshes the COVID-19 "Pais" sheet: country totals move up/down in the
# table (re-sorted by total cases), and the "updated at" timestamp changes.
#
# Because rows are kept sorted by total cases (column B) descending, a
# handful of countries swapped ranking positions between this refresh and
# the previous one. Rather than moving rows, every cell (country name +
# all of its numeric stats) for each affected row is rewritten in place
# with its new row's contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number -> full new row values (A:H) = Pais, Casos totales, Nuevos casos,
# Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes
$rowData = @{
    4   = @('Estados Unidos', 3068207, 28015, 1338696, 1595994, 0, 538, 133517)
    6   = @('India', 743481, 23135, 457045, 265783, 0, 479, 20653)
    10  = @('España', 299210, 341, 0, 0, 0, 4, 28392)
    19  = @('Alemania', 198310, 253, 182700, 6513, 0, 5, 9097)
    20  = @('Francia', 168810, 475, 77655, 61222, 0, 13, 29933)
    21  = @('Banglades', 168645, 3027, 78102, 88392, 0, 55, 2151)
    69  = @('Chequia', 12639, 73, 7910, 4378, 0, 1, 351)
    71  = @('Uzbekistan', 10587, 225, 6690, 3857, 0, 3, 40)
    91  = @('Costa Rica', 5486, 245, 1810, 3653, 0, 0, 23)
    93  = @('Mauritania', 5024, 76, 1944, 2945, 0, 2, 135)
    95  = @('Estado de Palestina', 4647, 306, 494, 4135, 0, 1, 18)
    108 = @('Paraguay', 2502, 46, 1193, 1289, 0, 0, 20)
    109 = @('Maldivas', 2501, 10, 2158, 331, 0, 0, 12)
    143 = @('Liberia', 917, 26, 394, 482, 0, 2, 41)
    144 = @('Montenegro', 907, 66, 320, 570, 0, 3, 17)
    149 = @('Santo Tome y Principe', 724, 3, 279, 432, 0, 0, 13)
    159 = @('Angola', 386, 40, 117, 248, 0, 2, 21)
    160 = @('Siria', 372, 0, 126, 232, 0, 0, 14)
    161 = @('Vietnam', 369, 0, 341, 28, 0, 0, 0)
    167 = @('Guyana', 278, 5, 121, 142, 0, 0, 15)
    209 = @('Islas Malvinas', 13, 0, 13, 0, 0, 0, 0)
    210 = @('Groenlandia', 13, 0, 13, 0, 0, 0, 0)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
}

# Update the "last refreshed" caption above the table.
$ws.Range("A1").Value = "Datos actualizados a 7 de Julio de 2020 a las 21:21"
